# Applies the "Updated symbol list on Wed Jan 25 09:29:44 UTC 2023 with GitHub Actions"
# edit: refreshed Price/Volume(1h) figures, plus a one-row insertion of "GateToken"
# at row 7 that shifted the KuCoinToken..LEO block down by one row (each row's
# Coin/Link moved to the next row while Price/Volume(1h) were independently refreshed).
#
# All Price (column D) and Volume(1h) (column E) entries are stored as literal text
# (t="inlineStr" in the original OOXML) so that exact formatting (trailing zeros,
# "%" suffixes, etc.) survives - a bare numeric-looking assignment would make Excel
# coerce the string to a Double and silently reformat/round it. Prefixing the string
# with a single quote forces Excel to keep it as text, matching the source file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'301.52"
$ws.Range("E2").Value = "'-4.41%"

# Row 3
$ws.Range("D3").Value = "'35.18"
$ws.Range("E3").Value = "'-2.29%"

# Row 4
$ws.Range("D4").Value = "'5.063"
$ws.Range("E4").Value = "'-1.45%"

# Row 5
$ws.Range("D5").Value = "'0.07967"
$ws.Range("E5").Value = "'-2.21%"

# Row 6
$ws.Range("D6").Value = "'1.924"
$ws.Range("E6").Value = "'-9.60%"

# Row 7
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'4.050"
$ws.Range("E7").Value = "'-2.27%"

# Row 8
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'7.731"
$ws.Range("E8").Value = "'-3.74%"

# Row 9
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.911"
$ws.Range("E9").Value = "'1.35%"

# Row 10
$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D10").Value = "'0.9228"
$ws.Range("E10").Value = "'-0.75%"

# Row 11
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.1283"
$ws.Range("E11").Value = "'26.09%"

# Row 12
$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").Value = "'0.1850"
$ws.Range("E12").Value = "'-1.65%"

# Row 13
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "'0.09675"
$ws.Range("E13").Value = "'4.93%"

# Row 14
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03632"
$ws.Range("E14").Value = "'1.23%"

# Row 15
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09842"
$ws.Range("E15").Value = "'-0.71%"

# Row 16
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001387"
$ws.Range("E16").Value = "'-4.13%"

# Row 17
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.005888"
$ws.Range("E17").Value = "'1.81%"

# Row 18
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.505"
$ws.Range("E18").Value = "'0.80%"

# Row 19
$ws.Range("D19").Value = "'0.3460"
$ws.Range("E19").Value = "'2.66%"

# Row 20
$ws.Range("D20").Value = "'0.1312"
$ws.Range("E20").Value = "'-1.51%"

# Row 21
$ws.Range("D21").Value = "'5.040"
$ws.Range("E21").Value = "'-2.38%"

# Row 22
$ws.Range("D22").Value = "'0.2405"
$ws.Range("E22").Value = "'9.57%"

# Row 23
$ws.Range("D23").Value = "'0.04529"
$ws.Range("E23").Value = "'-1.39%"

# Row 24
$ws.Range("D24").Value = "'0.001219"
$ws.Range("E24").Value = "'-2.50%"

# Row 25
$ws.Range("D25").Value = "'0.004808"
$ws.Range("E25").Value = "'1.73%"

# Row 26
$ws.Range("D26").Value = "'0.0001252"
$ws.Range("E26").Value = "'-0.17%"

# Row 27
$ws.Range("D27").Value = "'0.0003006"
$ws.Range("E27").Value = "'-33.40%"

# Row 39
$ws.Range("D39").Value = "'0.01902"
$ws.Range("E39").Value = "'-4.17%"

# Row 40
$ws.Range("D40").Value = "'0.04678"
$ws.Range("E40").Value = "'-4.82%"

# Row 41
$ws.Range("D41").Value = "'0.007535"
$ws.Range("E41").Value = "'-3.77%"

# Row 42
$ws.Range("D42").Value = "'0.009656"
$ws.Range("E42").Value = "'23.35%"

# Row 43
$ws.Range("D43").Value = "'0.1321"
$ws.Range("E43").Value = "'-5.56%"

# Row 44
$ws.Range("D44").Value = "'0.002113"
$ws.Range("E44").Value = "'-0.19%"

# Row 45
$ws.Range("D45").Value = "'0.01083"
$ws.Range("E45").Value = "'-6.99%"

# Row 46
$ws.Range("D46").Value = "'0.00006247"
$ws.Range("E46").Value = "'-4.29%"

# Row 47
$ws.Range("E47").Value = "'-0.18%"

# Row 48
$ws.Range("E48").Value = "'64.89%"

# Row 49
$ws.Range("E49").Value = "'-21.85%"

# Row 50
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'-0.18%"

# Row 51
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'-0.18%"
